$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A2").Value = "Última actualización: 01:59:52"
$ws1.Range("A3").Value = "Total filas: 4"

$ws1.Range("A6").Value = "01:59:52"
$ws1.Range("B6").Value = "01:59"
$ws1.Range("D6").Value = 0

$ws1.Range("A7").Value = "01:59:52"
$ws1.Range("D7").Value = 58

$ws1.Range("A8").Value = "01:59:52"
$ws1.Range("B8").Value = "03:02"
$ws1.Range("D8").Value = 63

$ws1.Range("A9").Value = "01:59:52"
$ws1.Range("B9").Value = "03:48"
$ws1.Range("C9").Value = "14_ABASTO"
$ws1.Range("D9").Value = 109
$ws1.Range("E9").Value = "LP1912"

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A2").Value = "Última actualización: 01:59:52"

$ws2.Range("A6").Value = "01:59:52"
$ws2.Range("D6").Value = 58

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A2").Value = "Última actualización: 01:59:52"
